$wb = $excel.ActiveWorkbook

# --- weibull sheet ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.04690973252514
$ws.Range("C2").Value = 0.0737316968516354
$ws.Range("B3").Value = -0.171946395700535
$ws.Range("C3").Value = 0.0388523711236188

# --- lognormal sheet ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.13237362326564
$ws.Range("C2").Value = 0.0882260712790033
$ws.Range("B3").Value = -0.691559258068776
$ws.Range("C3").Value = 0.0375498739940975

# --- llogis sheet ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.64016821659955
$ws.Range("C2").Value = 0.0855002121713666
$ws.Range("B3").Value = 1.1820943127845
$ws.Range("C3").Value = 0.0656643713394952

# --- gompertz sheet ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -1.86178907503081
$ws.Range("C2").Value = 0.0914459151885326
$ws.Range("B3").Value = -0.0646415225040878
$ws.Range("C3").Value = 0.0128324640774315

# --- weibull cov sheet ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.00543636312062147
$ws.Range("B2").Value = -0.00158773213442984
$ws.Range("A3").Value = -0.00158773213442984
$ws.Range("B3").Value = 0.00150950674192741

# --- lognormal cov sheet ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.00778383965332777
$ws.Range("B2").Value = -0.00251196077708807
$ws.Range("A3").Value = -0.00251196077708807
$ws.Range("B3").Value = 0.0014099930369726

# --- llogis cov sheet ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0073102862813487
$ws.Range("B2").Value = -0.000502977650851661
$ws.Range("A3").Value = -0.000502977650851661
$ws.Range("B3").Value = 0.00431180966341111

# --- gompertz cov sheet ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0083623554046683
$ws.Range("B2").Value = -0.00089930462489449
$ws.Range("A3").Value = -0.00089930462489449
$ws.Range("B3").Value = 0.000164672134298569
